$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "61.654.49"
$ws.Cells.Item(2, 5).Value = "  +5.34%  "

$ws.Cells.Item(3, 4).Value = "3.056.15"
$ws.Cells.Item(3, 5).Value = "  +2.75%  "

$ws.Cells.Item(4, 5).Value = "  -0.12%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "576.94"
$ws.Cells.Item(5, 5).Value = "  +3.28%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "142.60"
$ws.Cells.Item(6, 5).Value = "  +4.62%  "

$ws.Cells.Item(7, 5).Value = "  -0.22%  "

$ws.Cells.Item(8, 4).Value = "3.049.01"
$ws.Cells.Item(8, 5).Value = "  +2.69%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.524"
$ws.Cells.Item(9, 5).Value = "  +1.84%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.138"
$ws.Cells.Item(10, 5).Value = "  +5.81%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "5.45"
$ws.Cells.Item(11, 5).Value = "  +13.41%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.463"
$ws.Cells.Item(12, 5).Value = "  +2.28%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.0000237"
$ws.Cells.Item(13, 5).Value = "  +4.42%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "34.64"
$ws.Cells.Item(14, 5).Value = "  +3.41%  "

$ws.Cells.Item(15, 5).Value = "  -0.06%  "

$ws.Cells.Item(16, 4).Value = "3.559.46"
$ws.Cells.Item(16, 5).Value = "  +2.41%  "

$ws.Cells.Item(17, 5).Value = "  +3.63%  "

$ws.Cells.Item(18, 4).Value = "3.047.88"
$ws.Cells.Item(18, 5).Value = "  +2.23%  "

$ws.Cells.Item(19, 4).Value = "61.467.78"
$ws.Cells.Item(19, 5).Value = "  +4.83%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "447.59"
$ws.Cells.Item(20, 5).Value = "  +5.86%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "13.85"
$ws.Cells.Item(21, 5).Value = "  +2.79%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.730"
$ws.Cells.Item(22, 5).Value = "  +3.06%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "7.28"
$ws.Cells.Item(23, 5).Value = "  +2.62%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "13.59"
$ws.Cells.Item(24, 5).Value = "  +1.71%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "81.65"
$ws.Cells.Item(25, 5).Value = "  +2.13%  "

$ws.Cells.Item(26, 5).Value = "  +0.10%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "2.27"
$ws.Cells.Item(27, 5).Value = "  +8.44%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "0.998"
$ws.Cells.Item(28, 5).Value = "  -0.36%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "2.63"
$ws.Cells.Item(29, 5).Value = "  +4.44%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "8.04"
$ws.Cells.Item(30, 5).Value = "  +4.16%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "6.43"
$ws.Cells.Item(31, 5).Value = "  +6.37%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "26.45"
$ws.Cells.Item(32, 5).Value = "  +3.33%  "

$ws.Cells.Item(33, 5).Value = "  +6.70%  "

$ws.Cells.Item(34, 5).Value = "  +7.03%  "

$ws.Cells.Item(35, 5).Value = "  +3.50%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "6.05"
$ws.Cells.Item(36, 5).Value = "  +5.90%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "2.17"
$ws.Cells.Item(37, 5).Value = "  +5.21%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "50.04"
$ws.Cells.Item(38, 5).Value = "  +3.31%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "2.94"
$ws.Cells.Item(39, 5).Value = "  +7.78%  "

$ws.Cells.Item(40, 5).Value = "  +1.10%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "412.64"
$ws.Cells.Item(41, 5).Value = "  +3.11%  "

$ws.Cells.Item(42, 5).Value = "  +4.78%  "

$ws.Cells.Item(43, 4).Value = "2.782.32"
$ws.Cells.Item(43, 5).Value = "  +1.88%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.108"
$ws.Cells.Item(44, 5).Value = "  +0.51%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.262"
$ws.Cells.Item(45, 5).Value = "  +8.50%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "37.77"
$ws.Cells.Item(46, 5).Value = "  +19.01%  "

$ws.Cells.Item(47, 2).Value = "Fetch.AI"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "2.10"
$ws.Cells.Item(47, 5).Value = "  +4.75%  "

$ws.Cells.Item(48, 2).Value = "USDe"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.999"
$ws.Cells.Item(48, 5).Value = "  -0.04%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "123.84"
$ws.Cells.Item(49, 5).Value = "  -0.97%  "

$ws.Cells.Item(50, 5).Value = "  +1.93%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "23.99"
$ws.Cells.Item(51, 5).Value = "  +3.66%  "
